# Update cryptocurrency price (D) and volume-change (E) columns
# with the latest scraped values. Prices are forced to remain plain
# text (matching the workbook's existing inline-string convention)
# even when the new value looks numeric, so values like '26.50' or
# '0.0000164' don't get silently coerced into floating point numbers
# (which would drop trailing zeros / switch to scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($address, $text) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" '66.124.66'
Set-TextCell "E2" '  -0.72%  '
Set-TextCell "D3" '3.295.52'
Set-TextCell "E3" '  -0.81%  '
Set-TextCell "E4" '  +0.01%  '
Set-TextCell "D5" '585.77'
Set-TextCell "E5" '  +2.21%  '
Set-TextCell "D6" '181.22'
Set-TextCell "E6" '  -0.73%  '
Set-TextCell "D7" '0.649'
Set-TextCell "E7" '  +7.98%  '
Set-TextCell "E8" '  +0.02%  '
Set-TextCell "D9" '0.126'
Set-TextCell "E9" '  -3.13%  '
Set-TextCell "D10" '6.76'
Set-TextCell "E10" '  +1.69%  '
Set-TextCell "E11" '  +0.50%  '
Set-TextCell "D12" '3.866.77'
Set-TextCell "E12" '  -0.86%  '
Set-TextCell "E13" '  -4.57%  '
Set-TextCell "D14" '66.158.25'
Set-TextCell "D15" '26.50'
Set-TextCell "E15" '  -2.40%  '
Set-TextCell "D16" '0.0000164'
Set-TextCell "E16" '  -2.22%  '
Set-TextCell "D17" '3.285.70'
Set-TextCell "E17" '  -1.42%  '
Set-TextCell "D18" '433.33'
Set-TextCell "E18" '  -1.45%  '
Set-TextCell "D19" '13.28'
Set-TextCell "E19" '  -3.41%  '
Set-TextCell "D20" '5.51'
Set-TextCell "E20" '  -3.13%  '
Set-TextCell "D21" '7.43'
Set-TextCell "E21" '  -3.00%  '
Set-TextCell "D22" '72.35'
Set-TextCell "E22" '  -1.98%  '
Set-TextCell "E23" '  +0.10%  '
Set-TextCell "E24" '  +0.31%  '
Set-TextCell "D25" '3.434.03'
Set-TextCell "E25" '  -0.79%  '
Set-TextCell "D26" '0.512'
Set-TextCell "E26" '  -0.85%  '
Set-TextCell "E27" '  -3.96%  '
Set-TextCell "E28" '  +2.09%  '
Set-TextCell "D29" '8.88'
Set-TextCell "E29" '  -2.01%  '
Set-TextCell "D30" '0.999'
Set-TextCell "E30" '  -0.37%  '
Set-TextCell "E31" '  +0.56%  '
Set-TextCell "D32" '22.39'
Set-TextCell "E32" '  -2.31%  '
Set-TextCell "E33" '  -0.02%  '
Set-TextCell "E34" '  -3.00%  '
Set-TextCell "D35" '6.63'
Set-TextCell "E35" '  -2.37%  '
Set-TextCell "D36" '1.20'
Set-TextCell "E36" '  -2.59%  '
Set-TextCell "D37" '158.78'
Set-TextCell "E37" '  -0.81%  '
Set-TextCell "E38" '  -5.04%  '
Set-TextCell "D39" '26.61'
Set-TextCell "E39" '  -3.13%  '
Set-TextCell "E40" '  -3.61%  '
Set-TextCell "D41" '2.786.43'
Set-TextCell "E41" '  -0.98%  '
Set-TextCell "D42" '0.772'
Set-TextCell "E42" '  -2.35%  '
Set-TextCell "D43" '4.35'
Set-TextCell "E43" '  -2.78%  '
Set-TextCell "D44" '40.15'
Set-TextCell "E44" '  -0.13%  '
Set-TextCell "D45" '6.03'
Set-TextCell "E45" '  -3.30%  '
Set-TextCell "D46" '0.0661'
Set-TextCell "E46" '  -2.55%  '
Set-TextCell "E47" '  -1.14%  '
Set-TextCell "D48" '23.31'
Set-TextCell "E48" '  -3.91%  '
Set-TextCell "D49" '316.68'
Set-TextCell "E49" '  -0.73%  '
Set-TextCell "E50" '  -1.80%  '
Set-TextCell "E51" '  +5.64%  '
